$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "ksallmx"
$ws.Range("B14").Value = "bcnxd'lkj["
$ws.Range("B15").Value = "d"

$ws.Range("B15").Select()
